$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy style from existing header cell (H1) to new headers so formatting matches
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A1").Select()

# Data values for columns I and J, rows 2-12
$data = @(
    @(6, 7),
    @(8, 8),
    @(4, 5),
    @(10, 10),
    @(4, 6),
    @(6, 8),
    @(7, 8),
    @(6, 9),
    @(5, 7),
    @(4, 6),
    @(1, 2)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
